$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.244.51"
$ws.Range("E2").Value = "  -0.68%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.842.34"
$ws.Range("E3").Value = "  -0.44%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.75"
$ws.Range("E5").Value = "  +0.22%  "

# Row 6 - XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6282"
$ws.Range("E6").Value = "  -0.15%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.10%  "

# Row 8 - Dogecoin
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07439"
$ws.Range("E8").Value = "  -2.83%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2888"
$ws.Range("E9").Value = "  -1.04%  "

# Row 10 - Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.23"
$ws.Range("E10").Value = "  -2.26%  "

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07726"
$ws.Range("E11").Value = "  -0.14%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.841.38"
$ws.Range("E12").Value = "  -2.48%  "

# Row 13 - Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.979"
$ws.Range("E13").Value = "  -1.08%  "

# Row 14 - Polygon
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6759"
$ws.Range("E14").Value = "  -0.74%  "

# Row 15 - ShibaInu
$ws.Range("E15").Value = "  -5.17%  "

# Row 16 - Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.98"
$ws.Range("E16").Value = "  -1.87%  "

# Row 17 - Uniswap
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.125"
$ws.Range("E17").Value = "  -1.18%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "29.277.89"
$ws.Range("E18").Value = "  -1.08%  "

# Row 19 - BitcoinCash
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "227.38"
$ws.Range("E19").Value = "  -0.66%  "

# Row 20 - Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.25"
$ws.Range("E20").Value = "  -0.69%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.15%  "

# Row 22 - Chainlink
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.362"
$ws.Range("E22").Value = "  -1.40%  "

# Row 23 - BinanceUSD
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.02%  "

# Row 24 - Monero
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.60"
$ws.Range("E24").Value = "  +0.84%  "

# Row 25 - Stellar
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1371"
$ws.Range("E25").Value = "  -1.05%  "

# Row 26 - Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.375"
$ws.Range("E26").Value = "  -0.60%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.52"
$ws.Range("E27").Value = "  -1.20%  "

# Row 28 - Row28 (Toncoin->Hedera)
$ws.Range("B28").Value = "Hedera"
$ws.Range("C28").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06277"
$ws.Range("E28").Value = "  +11.79%  "

# Row 29 - Row29 (Hedera->Toncoin)
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.391"
$ws.Range("E29").Value = "  +0.81%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.49%  "

# Row 31 - Filecoin
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.076"
$ws.Range("E31").Value = "  -1.39%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.039"
$ws.Range("E32").Value = "  -0.63%  "

# Row 33 - LidoDAOToken
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.816"
$ws.Range("E33").Value = "  -1.47%  "

# Row 34 - ARBITRUM
$ws.Range("E34").Value = "  -2.49%  "

# Row 35 - ImmutableX
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6936"
$ws.Range("E35").Value = "  -0.86%  "

# Row 36 - HuobiToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.588"
$ws.Range("E36").Value = "  -0.14%  "

# Row 37 - MXToken
$ws.Range("E37").Value = "  +2.99%  "

# Row 38 - Maker
$ws.Range("D38").Value = "1.249.11"
$ws.Range("E38").Value = "  +1.60%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +0.38%  "

# Row 40 - FraxShare
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.510"
$ws.Range("E40").Value = "  +0.58%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9071"
$ws.Range("E41").Value = "  -0.13%  "

# Row 42 - PaxDollar
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9990"
$ws.Range("E42").Value = "  -0.13%  "

# Row 43 - RocketPoolETH
$ws.Range("D43").Value = "1.999.53"
$ws.Range("E43").Value = "  -17.21%  "

# Row 44 - Quant
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.30"
$ws.Range("E44").Value = "  -0.68%  "

# Row 45 - Aave
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.00"
$ws.Range("E45").Value = "  -0.03%  "

# Row 46 - Aptos
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.036"
$ws.Range("E46").Value = "  -2.45%  "

# Row 47 - Row47 (BabyDogeCoin->Algorand)
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1160"
$ws.Range("E47").Value = "  +0.42%  "

# Row 48 - Row48 (Algorand->BabyDogeCoin)
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000116"
$ws.Range("E48").Value = "  -1.49%  "

# Row 49 - EnergySwap
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.007"
$ws.Range("E49").Value = "  -0.18%  "

# Row 50 - TheSandbox
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3928"
$ws.Range("E50").Value = "  -2.55%  "

# Row 51 - RenderToken
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.655"
$ws.Range("E51").Value = "  -1.43%  "
